# Generate Report for Handoff
# Replace every occurrence of the old generated-id "8c8a8d3c-dd90-49cb-93c5-d580d206f868"
# with the new generated-id "9370d91d-79c1-4bf4-94db-e856dd5822b6" across the three
# report sheets (Overview, zh-cn, de-de), and bump the associated timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "8c8a8d3c-dd90-49cb-93c5-d580d206f868"
$newId = "9370d91d-79c1-4bf4-94db-e856dd5822b6"

# Content hash embedded in the generated .xlf handoff file names (shared by
# both the zh-cn and de-de rows). Old value was 2096a67755c6dc9b2fb2ef26cdca7b98a89edd51.
$newXlfHash = "e1f6217312596771566fb2880e8fffe16b4f843a"

# The hyperlinks on each sheet all point at the same (unchanged) external GitHub
# address -- only the display text needs to reflect the new id.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d1c2fc3870511fbb37391f15763af69afbddc16/e2e/$oldId.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-16 12:57:01"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newId.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-16 12:56:54"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newId.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.$newXlfHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-16 12:57:01"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newId.md")
